$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.639.31"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "3.718.70"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "656.27"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.423"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "3.719.85"
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +5.67%  "
$ws.Range("D15").Value = "4.411.17"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000269"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "96.519.79"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +16.78%  "
$ws.Range("D19").Value = "3.729.52"
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("E20").Value = "  +4.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.530"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "522.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000204"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.168"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.40%  "
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  +11.69%  "
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "662.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.602"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.162"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +23.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.979"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.86%  "
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.451"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("E51").Value = "  +1.76%  "
